$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-54 down to 19-55
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly record
$newDate = Get-Date -Year 2021 -Month 8 -Day 6 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = $newDate
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100112005
$ws.Cells.Item(18, 7).Value = "Puerro"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 180
$ws.Cells.Item(18, 11).Value = 7500
$ws.Cells.Item(18, 12).Value = 8000
$ws.Cells.Item(18, 13).Value = 7750
$ws.Cells.Item(18, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(18, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(18, 16).Value = 388
$ws.Cells.Item(18, 17).Value = 20
$ws.Cells.Item(18, 18).Value = "Hortaliza"
